$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3162
$ws.Range("K3").Value = 3140
$ws.Range("K4").Value = 645
$ws.Range("K5").Value = 207
$ws.Range("K6").Value = 3718
$ws.Range("K7").Value = 10872

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 207
$ws.Range("K3").Value = 221
$ws.Range("K6").Value = 235
$ws.Range("K7").Value = 721

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 119
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 427

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 128
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 115
$ws.Range("K7").Value = 375

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 10
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K5").Value = 19
$ws.Range("K8").Value = 721
$ws.Range("K10").Value = 63
$ws.Range("K11").Value = 232
$ws.Range("K12").Value = 17
$ws.Range("K14").Value = 57
$ws.Range("K16").Value = 35
$ws.Range("K19").Value = 329
$ws.Range("K20").Value = 248
$ws.Range("K23").Value = 101
$ws.Range("K26").Value = 18
$ws.Range("K27").Value = 111
$ws.Range("K29").Value = 566
$ws.Range("K30").Value = 35
$ws.Range("K32").Value = 17
$ws.Range("K33").Value = 427
$ws.Range("K34").Value = 52
$ws.Range("K37").Value = 375
$ws.Range("K41").Value = 95
$ws.Range("K42").Value = 375
$ws.Range("K43").Value = 96
$ws.Range("K44").Value = 104
$ws.Range("K46").Value = 22
$ws.Range("K48").Value = 133
$ws.Range("K49").Value = 66
$ws.Range("K51").Value = 122
$ws.Range("K52").Value = 298
$ws.Range("K54").Value = 211
$ws.Range("K55").Value = 115
$ws.Range("K57").Value = 35
$ws.Range("K62").Value = 4
$ws.Range("K63").Value = 38
$ws.Range("K65").Value = 255
$ws.Range("K66").Value = 40
$ws.Range("K67").Value = 430
$ws.Range("K73").Value = 99
$ws.Range("K76").Value = 165
$ws.Range("K77").Value = 77
$ws.Range("K78").Value = 140
$ws.Range("K79").Value = 279
$ws.Range("K80").Value = 35
$ws.Range("K85").Value = 513
$ws.Range("K86").Value = 71
$ws.Range("K88").Value = 128
$ws.Range("K89").Value = 142
$ws.Range("K91").Value = 114
$ws.Range("K92").Value = 43
$ws.Range("K93").Value = 42
$ws.Range("K94").Value = 137
$ws.Range("K95").Value = 178
$ws.Range("K97").Value = 96
$ws.Range("K99").Value = 193
$ws.Range("K101").Value = 10872

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 134
$ws.Range("K3").Value = 141
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 430

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K2").Value = 10
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 211

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 153
$ws.Range("K3").Value = 198
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 566

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 111
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 329

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 97
$ws.Range("K3").Value = 118
$ws.Range("K7").Value = 375

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 41
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 97
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 84
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K3").Value = 4
$ws.Range("K6").Value = 18

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 7
$ws.Range("K6").Value = 40

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 64
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 31
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 29
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 26
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 32
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 185
$ws.Range("K3").Value = 177
$ws.Range("K7").Value = 513

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 4
